$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename / add defined names ---------------------------------------
# Existing names get namespaced prefixes; two new ones are introduced for
# the new load-factor input and revenue output.
$wb.Names.Item("capacity").Name = "input.capacity"
$wb.Names.Item("output").Name = "output.energy"
$wb.Names.Add('input.load_factor', '=Sheet1!$B$3')
$wb.Names.Add('output.revenue', '=Sheet1!$B$8')

# --- New rows: Revenue input + Annual revenue output -------------------
$ws.Range("A7").Value = "Revenue"
$ws.Range("B7").Value = 50
$ws.Range("C7").Value = "£/MWh"

$ws.Range("A8").Value = "Annual revenue"
$ws.Range("B8").Formula = "=B7*output.energy"
$ws.Range("C8").Value = "£/yr"

# --- Cosmetic: widen column A, move the active selection ---------------
$ws.Columns("A").ColumnWidth = 15

$ws.Range("D17").Select()
